$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - this is the anchor for the
# block that needs to be removed (a trailing "scraped site footer" block):
#   <empty paragraph>
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: ... Creative Commons Attribution"
# The empty paragraph right before it, and the copyright paragraph right
# after it, are part of the same block and go away too. The blank
# paragraph that originally followed the copyright line (and the
# following page-break paragraph) must be left untouched.

$jupiterPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Ver no Jupiter")) {
        $jupiterPara = $p
        break
    }
}

$blankBefore = $jupiterPara.Previous(1)
$copyrightPara = $jupiterPara.Next(1)

$r = $d.Range($blankBefore.Range.Start, $copyrightPara.Range.End)
$r.Delete()
